# Update for release to deploy 0.1.1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Bump the Version value (row 3, column B)
$ws.Range("B3").Value = "0.1.1"

# Bump the Date value (row 8, column B)
$ws.Range("B8").Value = "2024-11-11T17:53:38-06:00"

# Insert a new metadata row above the existing "Description" row (row 11)
# so it becomes row 12, shifting everything below it down by one.
$ws.Rows("11:11").Insert()

# Copy the formatting of the (now shifted) row below onto the newly
# inserted blank row so it keeps the same body style/border as the
# rest of the table.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new "Jurisdiction" property row (value left blank).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
